$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.779.29'
$ws.Range("E2").Value = '  -1.78%  '

$ws.Range("D3").Value = '3.403.01'
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.402.65'
$ws.Range("E8").Value = '  -1.50%  '

$ws.Range("E9").Value = '  -6.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.62%  '

$ws.Range("E11").Value = '  -4.15%  '

$ws.Range("E12").Value = '  -4.94%  '

$ws.Range("D13").Value = '3.993.41'
$ws.Range("E13").Value = '  -1.43%  '

$ws.Range("E14").Value = '  -0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.08%  '

$ws.Range("E16").Value = '  -9.57%  '

$ws.Range("D17").Value = '63.929.27'
$ws.Range("E17").Value = '  -1.70%  '

$ws.Range("D18").Value = '3.420.75'
$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("E19").Value = '  -4.99%  '

$ws.Range("E20").Value = '  -4.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.62%  '

$ws.Range("E22").Value = '  -4.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.518'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.74%  '

$ws.Range("E26").Value = '  -2.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.24%  '

$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.96%  '

$ws.Range("E32").Value = '  -3.38%  '

$ws.Range("E33").Value = '  -2.98%  '

$ws.Range("E34").Value = '  -5.01%  '

$ws.Range("E35").Value = '  -4.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.816'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.63%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0727'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.11%  '

$ws.Range("D41").Value = '2.776.52'
$ws.Range("E41").Value = '  -3.52%  '

$ws.Range("E42").Value = '  -1.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.95%  '

$ws.Range("E46").Value = '  -4.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '325.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.91%  '

$ws.Range("E49").Value = '  -5.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.97%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.102'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.71%  '
